$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GAM Outputs")

# --- Row 7 ---
$ws.Range("A7").Value = 44561
$ws.Range("E7").Value = 0.55300000000000005
$ws.Range("F7").Value = 5112.8999999999996
$ws.Range("G7").Value = 10212.049999999999
$ws.Range("H7").Value = 26.177
$ws.Range("I7").Value = 5.2549999999999999

# --- Row 8 ---
$ws.Range("A8").Value = 44561
$ws.Range("E8").Value = 0.57299999999999995
$ws.Range("F8").Value = 5084.2
$ws.Range("G8").Value = 10151.33
$ws.Range("H8").Value = 25.106000000000002
$ws.Range("I8").Value = 5.4320000000000004
$ws.Range("P8").Value = 5.9859999999999998

# --- Row 9 ---
$ws.Range("A9").Value = 44561
$ws.Range("E9").Value = 0.57199999999999995
$ws.Range("F9").Value = 5090.3999999999996
$ws.Range("G9").Value = 10158.700000000001
$ws.Range("H9").Value = 25.762
$ws.Range("I9").Value = 4.71
$ws.Range("Q9").Value = 7.5709999999999997

# --- Row 10 ---
$ws.Range("A10").Value = 44561
$ws.Range("E10").Value = 0.59599999999999997
$ws.Range("F10").Value = 5052.1000000000004
$ws.Range("G10").Value = 10077.99
$ws.Range("H10").Value = 24.584
$ws.Range("I10").Value = 4.9509999999999996
$ws.Range("P10").Value = 5.7679999999999998
$ws.Range("Q10").Value = 7.7279999999999998

# --- Row 11 ---
$ws.Range("A11").Value = 44561
$ws.Range("E11").Value = 0.61799999999999999
$ws.Range("F11").Value = 5041.6000000000004
$ws.Range("G11").Value = 10019.280000000001
$ws.Range("H11").Value = 24.817
$ws.Range("I11").Value = 4.5030000000000001
$ws.Range("R11").Value = 25.076000000000001

# --- Sheet view / selection state ---
$ws.Activate() | Out-Null
$ws.Range("C13").Select() | Out-Null
